$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "LQFP32"

# --- Fill in newly-populated cells (Port / Pin# columns, plus EXTI rows) ---

# Row 12 (# = 11)
$ws.Range("B12").Value = "B"
$ws.Range("C12").Value = 7

# Row 13 (# = 12)
$ws.Range("B13").Value = "B"
$ws.Range("C13").Value = 6

# Row 14 (# = 13)
$ws.Range("B14").Value = "A"
$ws.Range("C14").Value = 12

# Row 15 (# = 14)
$ws.Range("B15").Value = "A"
$ws.Range("C15").Value = 3

# Row 16 (# = 15)
$ws.Range("B16").Value = "A"
$ws.Range("C16").Value = 4

# Row 17 (# = 16)
$ws.Range("B17").Value = "A"
$ws.Range("C17").Value = 6

# Row 18 (# = 17)
$ws.Range("B18").Value = "A"
$ws.Range("C18").Value = 7

# Row 19 (# = 18) - EXTI0
$ws.Range("B19").Value = "A"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = "EXTI0"

# Row 20 (# = 19) - EXTI1
$ws.Range("B20").Value = "A"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "EXTI1"

$ws.Range("I19").Value = "trigger on falling edge"

# Row 21 (# = 20) - EXTI2
$ws.Range("B21").Value = "A"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = "EXTI2"
$ws.Range("G21").Value = "input for emergency stop interrup ttrigger"

$ws.Range("I20").Value = "trigger on falling edge"
$ws.Range("I21").Value = "trigger on falling edge"

# Remove now-unused trailing rows 22-26
$ws.Range("A22:I26").EntireRow.Delete()

# Re-apply the autofilter over the used range
[void]$ws.Range("A1:I26").AutoFilter()
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=LQFP32!`$A`$1:`$I`$26")
$name.Visible = $false

# Restore selection state
[void]$ws.Range("I28").Select()
